$d = $word.ActiveDocument

# 1) Update the "Please note..." paragraph text to append the
#    "(Questions 1-5)" qualifier, then add a trailing bold run that
#    contains just a single space (matching the diff's extra <w:r>).
$d.Content.Find.Execute(
    "Please note that the steps show rounded numbers, but that the final answers to the problems are calculated without rounding.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Please note that the steps show rounded numbers, but that the final answers to the problems are calculated without rounding. (Questions 1-5)",
    2) | Out-Null

$notePara = $d.Paragraphs.Item(4)
$tail = $notePara.Range
$tail.Collapse(0)
$tail.InsertAfter(" ")
# Force the inserted space to live in its own run (not merge back into
# the preceding run) by toggling Bold off then on.
$tail.Bold = 0
$tail.Bold = 1

# 2) Remove the now-duplicate "1 | Mode | The most frequently occurring
#    value" row from the table (first table in the document).
$t = $d.Tables.Item(1)
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $row = $t.Rows.Item($i)
    $partText = $row.Cells.Item(2).Range.Text
    if ($partText -like "Mode*") {
        $row.Delete()
        break
    }
}

# 3) Strip the erroneous trailing " Mode -XX.XXXXXX" figure from each
#    company's solution cell (Questions 2a-2e).
$d.Content.Find.Execute("Mean: 21.276 Median: 13.433 Mode -48.837209", $true, $false, $false, $false, $false, $true, 1, $false, "Mean: 21.276 Median: 13.433", 2) | Out-Null
$d.Content.Find.Execute("Mean: 33.482 Median: 20.838 Mode -62.837689", $true, $false, $false, $false, $false, $true, 1, $false, "Mean: 33.482 Median: 20.838", 2) | Out-Null
$d.Content.Find.Execute("Mean: 41.122 Median: 25.558 Mode -71.050584", $true, $false, $false, $false, $false, $true, 1, $false, "Mean: 41.122 Median: 25.558", 2) | Out-Null
$d.Content.Find.Execute("Mean: 0.706 Median: 1.892 Mode -44.416873", $true, $false, $false, $false, $false, $true, 1, $false, "Mean: 0.706 Median: 1.892", 2) | Out-Null
$d.Content.Find.Execute("Mean: -1.084 Median: -3.796 Mode -39.686099", $true, $false, $false, $false, $false, $true, 1, $false, "Mean: -1.084 Median: -3.796", 2) | Out-Null
